$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.062.18'
$ws.Range('E2').Value = '  +1.66%  '
$ws.Range('D3').Value = '3.899.71'
$ws.Range('E3').Value = '  +3.46%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'470.47"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +10.30%  '
$ws.Range('D6').Value = "'143.76"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.92%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = "'0.732"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.77%  '
$ws.Range('D10').Value = "'0.162"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.82%  '
$ws.Range('D11').Value = "'0.0000334"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +10.11%  '
$ws.Range('D12').Value = "'43.10"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.70%  '
$ws.Range('D13').Value = '4.514.30'
$ws.Range('E13').Value = '  +3.26%  '
$ws.Range('D14').Value = "'10.35"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.08%  '
$ws.Range('D15').Value = "'15.00"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.24%  '
$ws.Range('D16').Value = '3.892.72'
$ws.Range('E16').Value = '  +3.90%  '
$ws.Range('E17').Value = '  -0.29%  '
$ws.Range('D18').Value = "'19.79"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('E19').Value = '  +3.76%  '
$ws.Range('D20').Value = '67.255.95'
$ws.Range('E20').Value = '  +1.74%  '
$ws.Range('D21').Value = "'430.89"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.95%  '
$ws.Range('B22').Value = 'InternetComputer(DFINITY)'
$ws.Range('C22').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D22').Value = "'14.62"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.10%  '
$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D23').Value = "'3.36"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.68%  '
$ws.Range('D24').Value = "'87.92"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.44%  '
$ws.Range('D25').Value = "'3.57"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +9.33%  '
$ws.Range('E26').Value = '  +5.42%  '
$ws.Range('E27').Value = '  +5.93%  '
$ws.Range('D28').Value = "'10.01"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.19%  '
$ws.Range('D29').Value = "'9.57"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.17%  '
$ws.Range('D30').Value = "'726.15"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.39%  '
$ws.Range('D31').Value = "'13.70"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('E32').Value = '  +0.17%  '
$ws.Range('D33').Value = "'2.81"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.48%  '
$ws.Range('D34').Value = "'43.10"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.43%  '
$ws.Range('E35').Value = '  +4.47%  '
$ws.Range('D36').Value = "'57.10"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.71%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('E38').Value = '  +20.52%  '
$ws.Range('E39').Value = '  -4.76%  '
$ws.Range('D40').Value = "'0.0474"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.39%  '
$ws.Range('D41').Value = "'3.07"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.43%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('E43').Value = '  +5.04%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').Value = "'2.57"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.21%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').Value = "'1.00"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('E46').Value = '  +4.22%  '
$ws.Range('D47').Value = "'2.16"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.88%  '
$ws.Range('E48').Value = '  +1.15%  '
$ws.Range('D49').Value = "'3.17"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.27%  '
$ws.Range('D50').Value = "'144.51"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.17%  '
$ws.Range('E51').Value = '  +3.79%  '
